$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing step/result text (rows 2 and 3 -> C2/D2, C3/D3)
$ws.Range("C2").Value = 'Step 1: while logged out try to enter the "my assessments" page'
$ws.Range("D2").Value = "I am returned to the login screen"

$ws.Range("C3").Value = "Step 2: Login as a user with the correct role"
$ws.Range("D3").Value = "I am redirected to the dashboard of the user"

$ws.Range("C4").Value = ' Step 3: Click "My Assessments" '
$ws.Range("D4").Value = "I am showed a list of assessments all based around me"

# Add new rows 5 and 6 with additional steps/results
$ws.Range("C5").Value = "Step 4: From the url view a list of someone elses assessments"
$ws.Range("D5").Value = "I am redirected to a list of my assessments(I should not be able to see a list of other people's assessments)"

$ws.Range("C6").Value = "Step 5: From the url change the myassessments id to nothing"
$ws.Range("D6").Value = "Redirected to a page saying access denied"

# Update selection to the C2:D4 range
$ws.Range("C2:D4").Select()
